# newly added iAuthor TC's
# Updates the generated candidate/marker credentials and candidate IDs
# on the "users" worksheet (rows 2 and 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Marker / "MR")
$ws.Range("A2").Value = 'QcymZ839'
$ws.Range("B2").Value = 23110950
$ws.Range("C2").Value = 'dbtjrpl67'
$ws.Range("D2").Value = 'yWn$2#B5'
$ws.Range("F2").Value = 'eJVKwVSb'
$ws.Range("G2").Value = 'wJiu'

# Row 3 (Candidate)
$ws.Range("A3").Value = 'Gyvdz314'
$ws.Range("B3").Value = 23110949
$ws.Range("C3").Value = 'utumbik60'
$ws.Range("D3").Value = 'Q&wj7B#2'
$ws.Range("F3").Value = 'bETtBKEq'
$ws.Range("G3").Value = 'eQHk'
